# Edit: Wed, Apr 29, 2020 5:07:12 PM
#
# 1) Three tables (on the slides that hold the "Google Shape" tables) get
#    their table style switched from the deck's single custom style
#    ({759EFD15-B15A-426B-8344-1CD0E4100226}) to the built-in style
#    {5318270B-0610-40B3-9564-4EBEE0E9146C}.
# 2) The presentation's theme colour scheme is swapped from the "Red
#    Violet"/Integral palette back to the standard "Office" palette.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Re-style every table in the deck that currently uses the custom
#    table style with the new style id.
# ---------------------------------------------------------------------
$oldStyleId = "{759EFD15-B15A-426B-8344-1CD0E4100226}"
$newStyleId = "{5318270B-0610-40B3-9564-4EBEE0E9146C}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId, $true)
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Restore the standard "Office" theme colours (dk1/lt1/dk2/lt2/
#    accent1-6/hlink/folHlink), replacing the current "Red Violet"
#    palette, via the slide's ThemeColorScheme (indices 1-12 map onto
#    dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink in order).
# ---------------------------------------------------------------------
function RGBValue([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officeColors = @(
    (RGBValue 0x00 0x00 0x00),  # dk1      000000
    (RGBValue 0xFF 0xFF 0xFF),  # lt1      FFFFFF
    (RGBValue 0x44 0x54 0x6A),  # dk2      44546A
    (RGBValue 0xE7 0xE6 0xE6),  # lt2      E7E6E6
    (RGBValue 0x5B 0x9B 0xD5),  # accent1  5B9BD5
    (RGBValue 0xED 0x7D 0x31),  # accent2  ED7D31
    (RGBValue 0xA5 0xA5 0xA5),  # accent3  A5A5A5
    (RGBValue 0xFF 0xC0 0x00),  # accent4  FFC000
    (RGBValue 0x44 0x72 0xC4),  # accent5  4472C4
    (RGBValue 0x70 0xAD 0x47),  # accent6  70AD47
    (RGBValue 0x05 0x63 0xC1),  # hlink    0563C1
    (RGBValue 0x95 0x4F 0x72)   # folHlink 954F72
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
